$d = $word.ActiveDocument

# --- Change 1: paragraph 2 ("I'm Amey Desai here, and I'm delighted ...")
# Remove the "Amey Desai here, and " lead-in, collapsing the run-split
# (and spell-check proofErr markup) that surrounded the name into a
# single run that starts directly with "I'm delighted ...".
$p2 = $d.Paragraphs(2)
$p2Start = $p2.Range.Start
$p2End = $p2.Range.End - 1
$target2 = "I'm delighted to provide you some information about your business. I appreciate you asking me the questions that are leading. It was helpful to see the kinds of insights you might anticipate from the data. I have every confidence that when you choose how to approach your next business possibilities, you will find the analysis to be convincing and helpful."
$d.Range($p2Start, $p2End).Text = $target2

# --- Change 2: paragraph 7 ("The second graph shows ...")
# No wording changes here -- just collapse the run split / proofErr
# markup around "realise" into a single run (identical text).
$p7 = $d.Paragraphs(7)
$p7Start = $p7.Range.Start
$p7End = $p7.Range.End - 1
$r7 = $d.Range($p7Start, $p7End)
$text7 = $r7.Text
$r7.Text = $text7 + "X"
$p7b = $d.Paragraphs(7)
$p7NewEnd = $p7b.Range.End - 1
$d.Range($p7NewEnd - 1, $p7NewEnd).Text = ""

# --- Change 3: add a signature name "Supin Hooda" in the final
# (previously empty) paragraph of the letter.
$n = $d.Paragraphs.Count
$last = $d.Paragraphs($n)
$rLast = $last.Range
$rLast.Text = "Supin Hooda"
$rLast.Font.Color = 2434341
